$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$shape = $s.Shapes.Item(2)

# Move the content placeholder down (EMU -> points: 1 pt = 12700 EMU)
$shape.Top = 721230 / 12700

# Replace the single-paragraph, multi-run URL text with three separate
# hyperlinked paragraphs.
$tr = $shape.TextFrame.TextRange
$tr.Text = "https://github.com/ariutta/svg-pan-zoom#svg-pan-zoom-library`r" + `
           "http://flask.palletsprojects.com/en/1.1.x/`r" + `
           "https://api.jquery.com/"

$tr.Paragraphs(1, 1).ActionSettings(1).Hyperlink.Address = "https://github.com/ariutta/svg-pan-zoom#svg-pan-zoom-library"
$tr.Paragraphs(2, 1).ActionSettings(1).Hyperlink.Address = "http://flask.palletsprojects.com/en/1.1.x/"
$tr.Paragraphs(3, 1).ActionSettings(1).Hyperlink.Address = "https://api.jquery.com/"
